# "add 1989 tv data"
# Append three new release rows for the 1989 (Taylor's Version) era to the
# `releases` sheet: two new album editions and one new track belonging to
# the existing "1989 (Taylor's Version)" album.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New album: "1989 (Taylor's Version) [Deluxe]" released 2023-10-27 (45226)
$ws.Cells.Item(133, 1).Value = "1989 (Taylor's Version) [Deluxe]"
$ws.Cells.Item(133, 4).Value = 45226

# New album: "1989 (Taylor's Version) [Tangerine Edition]" released 2023-10-27
$ws.Cells.Item(134, 1).Value = "1989 (Taylor's Version) [Tangerine Edition]"
$ws.Cells.Item(134, 4).Value = 45226

# New promotional track off the existing "1989 (Taylor's Version)" album
$ws.Cells.Item(135, 1).Value = "1989 (Taylor's Version)"
$ws.Cells.Item(135, 2).Value = "Slut! (Taylor's Version) [From The Vault]"
$ws.Cells.Item(135, 3).Value = 45226

# Widen column A to fit the longer album names now in the sheet
$ws.Columns.Item(1).ColumnWidth = 36.14

# Update the view to where the author's selection ended up after entry
$ws.Range("B136").Select() | Out-Null
